# Caso2 - Escenario IFrames
# Adds a second scenario sheet ("Scenario2") to the data-driven test workbook,
# duplicated from "Scenario1" (to preserve sheet-level formatting/page setup),
# with a new username/password data row, and updates both sheets' selections
# and active-tab state to reflect Scenario2 being the newly active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Record the selection left behind on Scenario1 before switching away from it.
$ws1.Range("A2").Select() | Out-Null

# Duplicate Scenario1 right after itself to become Scenario2.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Scenario2"

# New credentials row for the IFrames scenario (headers are inherited as-is).
$ws2.Range("A2").Value = "MartiQA"
$ws2.Range("B2").Value = "Pepe1234#"
$ws2.Range("C2").Value = "Test01"
$ws2.Range("D2").Value = "Test01"

# Leave the new sheet active with B3 selected.
$ws2.Range("B3").Select() | Out-Null
